$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop all the old data rows (2-8) - the fingerprinting/ROX rework only keeps
# the header plus two UMI adapter rows, so rebuild the body from scratch.
$ws.Rows("2:8").Delete()

# Row 2 - MatrixTube075 / UMIADAPTERU / 3 / 2 / Inline First Read
$ws.Range("A2").Value = "MatrixTube075"
$ws.Range("B2").Value = "UMIADAPTERU"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "Inline First Read"

# Row 3 - MatrixTube075 / UMIADAPTERU / 3 / 2 / Before Second Read
$ws.Range("A3").Value = "MatrixTube075"
$ws.Range("B3").Value = "UMIADAPTERU"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "Before Second Read"

# New highlighted font for the fingerprinting/ROX related cell E2.
$ws.Range("E2").Font.Name = "Helvetica Neue"
$ws.Range("E2").Font.Size = 12
$ws.Range("E2").Font.Color = 3355443

# Row 2 gets a taller custom row height to match the bigger font.
$ws.Rows(2).RowHeight = 15

# Update selection to reflect the new active cell.
$ws.Range("E2").Select()
